$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.487.19'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.645.69'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '302.59'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '0.3838'
$ws.Range("E7").Value = '  +0.66%  '
$ws.Range("D8").Value = '0.3594'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '51.01'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").Value = '0.08159'
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").Value = '1.227'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '22.28'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '6.430'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '7.435'
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '0.00001218'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '1.643.40'
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '97.47'
$ws.Range("E18").Value = '  +2.89%  '
$ws.Range("D19").Value = '0.07010'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '6.753'
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").Value = '17.51'
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '12.60'
$ws.Range("E23").Value = '  +1.61%  '
$ws.Range("D24").Value = '23.488.16'
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("D25").Value = '2.485'
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").Value = '3.025'
$ws.Range("E26").Value = '  -2.65%  '
$ws.Range("D27").Value = '21.17'
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("D28").Value = '152.82'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").Value = '5.233'
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = '133.88'
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").Value = '1.829.23'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '7.077'
$ws.Range("E32").Value = '  +9.02%  '
$ws.Range("D33").Value = '2.249'
$ws.Range("E33").Value = '  +4.96%  '
$ws.Range("D34").Value = '12.27'
$ws.Range("E34").Value = '  +6.66%  '
$ws.Range("D35").Value = '1.053'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("D36").Value = '0.02784'
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("D37").Value = '0.2498'
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").Value = '0.08765'
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("D39").Value = '6.051'
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D40").Value = '0.06971'
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '13.06'
$ws.Range("E41").Value = '  +8.01%  '
$ws.Range("D42").Value = '0.6963'
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("D43").Value = '1.334'
$ws.Range("E43").Value = '  +1.27%  '
$ws.Range("D44").Value = '15.87'
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("D45").Value = '0.6488'
$ws.Range("E45").Value = '  +1.75%  '
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = '2.289'
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("D48").Value = '3.954'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").Value = '0.07870'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = '127.83'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = '1.174'
$ws.Range("E51").Value = '  -0.46%  '
